# Update transition rule summary tables to include tri proximity tables.
$wb = $excel.ActiveWorkbook

# --- Means sheet ---
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("G2").Value = 66
$wsMeans.Range("G3").Value = 27
$wsMeans.Range("G4").Value = 7.1
$wsMeans.Range("G8").Value = 4.9
$wsMeans.Range("G9").Value = 80

# --- Standard Deviations sheet ---
$wsSd = $wb.Worksheets.Item("Standard Deviations")
$wsSd.Range("G5").Value = 6.2
$wsSd.Range("G7").Value = 7
$wsSd.Range("G8").Value = 8.1
